$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)

# Set explicit position/size override on the body placeholder (previously inherited from layout)
$sh.Left = 27.025630950927734
$sh.Top = 111.1197280883789
$sh.Width = 652.8189697265625
$sh.Height = 382.0394287109375

# Append a new bullet paragraph (same indent level as the last one) with the
# additional explanation text about unnecessary complexity with a single implementation
$tr = $sh.TextFrame.TextRange
[void]$tr.InsertAfter("`rBei nur einer Implementierung des Interfaces  unnötige Komplexität")
